$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting the existing data rows (2-42) down to (3-43)
$ws.Rows.Item(2).Insert()

# The newly inserted row inherits formatting from the row above it (the header),
# so copy the number formatting from row 3 (the row right below, which still has
# the original data-row styling) onto the new row 2.
$ws.Range("A3:C3").Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row 2 with the latest reported day's figures
$ws.Cells.Item(2, 1).Value = 43949
$ws.Cells.Item(2, 2).Value = 1635
$ws.Cells.Item(2, 3).Value = 165951
